# Update "想去人数" (want-to-go count) figures in column F across sheets.
# Sheet "展览" (Exhibitions) rows 4-42
$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F4").Value = 3666
$wsExhibition.Range("F5").Value = 3666
$wsExhibition.Range("F6").Value = 274
$wsExhibition.Range("F7").Value = 5192
$wsExhibition.Range("F8").Value = 552
$wsExhibition.Range("F9").Value = 384
$wsExhibition.Range("F11").Value = 708
$wsExhibition.Range("F13").Value = 111
$wsExhibition.Range("F16").Value = 326
$wsExhibition.Range("F21").Value = 365
$wsExhibition.Range("F22").Value = 4959
$wsExhibition.Range("F26").Value = 6084
$wsExhibition.Range("F29").Value = 3236
$wsExhibition.Range("F31").Value = 721
$wsExhibition.Range("F36").Value = 1069
$wsExhibition.Range("F37").Value = 87
$wsExhibition.Range("F40").Value = 890
$wsExhibition.Range("F41").Value = 1050
$wsExhibition.Range("F42").Value = 2039

# Sheet "本地生活" (Local Life) row 3
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F3").Value = 1131

# Sheet "全部类型" (All Types) rows 4, 7-48
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1131
$wsAll.Range("F7").Value = 3666
$wsAll.Range("F8").Value = 3666
$wsAll.Range("F9").Value = 274
$wsAll.Range("F10").Value = 5192
$wsAll.Range("F11").Value = 552
$wsAll.Range("F12").Value = 384
$wsAll.Range("F14").Value = 708
$wsAll.Range("F16").Value = 111
$wsAll.Range("F19").Value = 326
$wsAll.Range("F25").Value = 365
$wsAll.Range("F26").Value = 4959
$wsAll.Range("F30").Value = 6084
$wsAll.Range("F33").Value = 3236
$wsAll.Range("F35").Value = 721
$wsAll.Range("F41").Value = 1069
$wsAll.Range("F42").Value = 87
$wsAll.Range("F45").Value = 890
$wsAll.Range("F46").Value = 1050
$wsAll.Range("F48").Value = 2039
